$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (e.g. "555.71").
# Excel would silently coerce a bare numeric-looking string into a Number
# cell, but the source workbook keeps these as text. Prefixing with an
# apostrophe forces text entry; ClearFormats() then drops the resulting
# quotePrefix style so the cell formatting stays identical to the original.

$ws.Range("D2").Value = "62.233.19"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "2.424.78"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'555.71"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "'143.47"
$ws.Range("E6").Value = "  +5.37%  "
$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("D9").Value = "2.425.53"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "'5.40"
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").Value = "'26.32"
$ws.Range("E14").Value = "  +7.48%  "
$ws.Range("E15").Value = "  +9.64%  "
$ws.Range("D16").Value = "2.862.30"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "62.038.03"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").Value = "2.423.34"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").Value = "'11.12"
$ws.Range("E19").Value = "  +5.33%  "
$ws.Range("D20").Value = "'325.29"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").Value = "'6.76"
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("D25").Value = "'64.98"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("D26").Value = "'9.20"
$ws.Range("E26").Value = "  +11.45%  "
$ws.Range("D27").Value = "'573.10"
$ws.Range("E27").Value = "  +15.15%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.538.03"
$ws.Range("E28").Value = "  +3.45%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'8.38"
$ws.Range("E30").Value = "  +6.07%  "
$ws.Range("E31").Value = "  +9.64%  "
$ws.Range("E32").Value = "  +6.11%  "
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("E34").Value = "  +5.17%  "
$ws.Range("E35").Value = "  +4.69%  "
$ws.Range("D36").Value = "'5.76"
$ws.Range("E36").Value = "  +10.80%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.83"
$ws.Range("E38").Value = "  +5.93%  "
$ws.Range("D39").Value = "'1.94"
$ws.Range("E39").Value = "  +7.83%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").Value = "'146.99"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'41.59"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").Value = "'2.32"
$ws.Range("E45").Value = "  +12.43%  "
$ws.Range("D46").Value = "'151.68"
$ws.Range("E46").Value = "  +7.22%  "
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("E48").Value = "  +7.18%  "
$ws.Range("D49").Value = "'20.52"
$ws.Range("E49").Value = "  +8.32%  "
$ws.Range("E50").Value = "  +4.68%  "
$ws.Range("D51").Value = "'0.0914"
$ws.Range("E51").Value = "  +1.80%  "

# Strip the auto-applied "quote prefix" formatting picked up above so the
# cells end up with the same (default) style as before the edit.
$numericTextCells = @("D5","D6","D8","D12","D14","D19","D20","D21","D22","D25","D26","D27","D29","D30","D36","D37","D38","D39","D42","D44","D45","D46","D49","D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}
